# Generate Report for Handoff
# The 82d82f6b-1ec0-4ac5-a4f0-77bc32a55d74.md item moved from "Handed back: in
# sync with en-US" to "Ready for handoff": a newer handback arrived, got
# picked up by the handoff job, and since the live en-US source has since
# moved on again, the job also recorded a "version is not latest" error.

$wb = $excel.ActiveWorkbook

$status = "Ready for handoff"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e7514cb58b62a2929cfa2c739c0a50291ae34e05/e2e/82d82f6b-1ec0-4ac5-a4f0-77bc32a55d74.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8f0c68e631147bb42e8d27a334ae9a61c275aeee/e2e/82d82f6b-1ec0-4ac5-a4f0-77bc32a55d74.md."

# ---- Overview sheet: row 3 is the 82d82f6b-...md file ----
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E3").Value = $status
$ov.Range("F3").Value = $status
$ov.Range("G3").Value = "2016-08-28 16:48:57"

# ---- zh-cn sheet: row 3 ----
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C3").Value = $status
$zh.Range("H3").Value = "2016-08-28 16:48:53"
$zh.Range("P3").Value = $errorDetail
$zh.Columns.Item(16).ColumnWidth = 39.17

# ---- de-de sheet: row 3 ----
$de = $wb.Worksheets.Item("de-de")
$de.Range("C3").Value = $status
$de.Range("H3").Value = "2016-08-28 16:48:57"
$de.Range("P3").Value = $errorDetail
$de.Columns.Item(16).ColumnWidth = 39.17
